$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AMOGAS Data 1"
$ws.Range("B2").Value = "BELLO"
$ws.Range("C2").Value = "Sample Data 1"
$ws.Range("D2").Value = "Sample Data 1"
$ws.Range("E2").Value = "Sample Data 1"
$ws.Range("F2").Value = "sam17.bello@ymail.com"
$ws.Range("G2").Value = "Sample data 2"
$ws.Range("H2").Value = "Sample data 2"
$ws.Range("I2").Value = "Sample data 2"
$ws.Range("J2").Value = "Sample data 2"
$ws.Range("K2").Value = "PHP 300"
$ws.Range("L2").Value = "PHP 20000"
$ws.Range("M2").Value = "FINAL PHP 20000"
$ws.Range("N2").Value = "₱26502.00"
$ws.Range("O2").Value = "₱857.61"
$ws.Range("P2").Value = "₱1192.59"
$ws.Range("Q2").Value = "PHP 200321300"
$ws.Range("R2").Value = "₱88.34"
$ws.Range("S2").Value = "₱100.00"
$ws.Range("T2").Value = "₱2258.54"
$ws.Range("U2").Value = "₱24243.46"
$ws.Range("V2").Value = "Sample data 2"
$ws.Range("W2").Value = "₱857.61"
$ws.Range("X2").Value = "₱1192.59"
$ws.Range("Y2").Value = 0
